$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in column C (the "Förändrad" column).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

# All data rows (2 through last row) in column C currently hold 45203
# (2023-10-04) and need to be bumped to 45204 (2023-10-05).
$ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item($lastRow, 3)).Value = 45204
